# Add a new "canonical SMILES" column (D) to the microstate list worksheet,
# mirroring the "canonical isomeric SMILES" column (C) except for
# SM16_micro005, whose canonical (non-isomeric) SMILES differs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column D - copy value, then copy formatting from C2.
$ws.Range("D2").Value = "canonical SMILES"
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)  # xlPasteFormats

# Row 3: SM16_micro001
$ws.Cells.Item(3, 4).Value = "c1cc(c(c(c1)Cl)C(=[OH+])N=c2cc[nH]cc2)Cl"
$ws.Cells.Item(3, 3).Copy()
$ws.Cells.Item(3, 4).PasteSpecial(-4122)

# Row 4: SM16_micro002
$ws.Cells.Item(4, 4).Value = "c1cc(c(c(c1)Cl)C(=O)Nc2ccncc2)Cl"
$ws.Cells.Item(4, 3).Copy()
$ws.Cells.Item(4, 4).PasteSpecial(-4122)

# Row 5: SM16_micro003
$ws.Cells.Item(5, 4).Value = "c1cc(c(c(c1)Cl)C(=[OH+])[N-]c2ccncc2)Cl"
$ws.Cells.Item(5, 3).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4122)

# Row 6: SM16_micro004
$ws.Cells.Item(6, 4).Value = "c1cc(c(c(c1)Cl)C(=[NH+]c2cc[nH+]cc2)[O-])Cl"
$ws.Cells.Item(6, 3).Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4122)

# Row 7: SM16_micro005 -- canonical SMILES differs from the isomeric SMILES
$ws.Cells.Item(7, 4).Value = "c1cc(c(c(c1)Cl)C(=Nc2ccncc2)[O-])Cl"
$ws.Cells.Item(7, 3).Copy()
$ws.Cells.Item(7, 4).PasteSpecial(-4122)

# Row 8: SM16_micro006
$ws.Cells.Item(8, 4).Value = "c1cc(c(c(c1)Cl)C(=[OH+])Nc2ccncc2)Cl"
$ws.Cells.Item(8, 3).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4122)

# Row 9: SM16_micro007
$ws.Cells.Item(9, 4).Value = "c1cc(c(c(c1)Cl)C(=O)N=c2cc[nH]cc2)Cl"
$ws.Cells.Item(9, 3).Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4122)

# Row 10: SM16_micro008
$ws.Cells.Item(10, 4).Value = "c1cc(c(c(c1)Cl)C(=[OH+])Nc2cc[nH+]cc2)Cl"
$ws.Cells.Item(10, 3).Copy()
$ws.Cells.Item(10, 4).PasteSpecial(-4122)

# Match the new column width from the diff. Excel COM snaps column widths to
# a pixel grid, so we use the input value that rounds to the closest
# achievable width to the target 36.85546875 (36.833333333333336).
$ws.Columns.Item(4).ColumnWidth = 36.0
